# Update row 28 metrics (2025Q2) on Sheet1 with refreshed numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = 414
$ws.Range("D28").Value = 47
$ws.Range("E28").Value = 367
$ws.Range("F28").Value = 7.320872274143301
